# GPLIM-3541: add Material Type as required header for Manifest uploads
#
# Adds a new "Material Type" column (G) to the manifest worksheet:
#   - G1 header "Material Type" styled like the existing SAMPLE_TYPE header
#     (bold white text centered on a solid black fill)
#   - G2:G24 populated with "DNA:Genomic", centered, matching the existing
#     centered-data style already used by column C
#   - selection moved to the new column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 24

# Header cell: bold white font on black fill, centered - mirrors F1's look.
$header = $ws.Range("G1")
$header.Value = "Material Type"
$header.Font.Name = "MS Sans Serif"
$header.Font.Size = 10
$header.Font.Bold = $true
$header.Font.Color = 16777215
$header.Interior.Color = 0
$header.HorizontalAlignment = -4108

# Data cells: centered text, same look as column C's values.
$data = $ws.Range("G2:G$lastRow")
$data.Value = "DNA:Genomic"
$data.HorizontalAlignment = -4108

# Move the active selection onto the newly added column, matching what
# Excel leaves selected right after inserting/filling this column.
$ws.Range("G1:G$lastRow").Select()
